$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# Hunk 0: ALC row 135
$ws_ALC.Range("H135").Value = 75001960
$ws_ALC.Range("I135").Value = 31251720
$ws_ALC.Range("J135").Value = 250002900
$ws_ALC.Range("K135").Value = 281265480
$ws_ALC.Range("L135").Value = 2250026100
$ws_ALC.Range("M135").Value = -281262945
$ws_ALC.Range("N135").Value = -2250031170

# Hunk 1: ALC row 137
$ws_ALC.Range("H137").Value = 3176.8438
$ws_ALC.Range("I137").Value = 1885.3914
$ws_ALC.Range("J137").Value = 6477.222
$ws_ALC.Range("K137").Value = 5656.174199999999
$ws_ALC.Range("L137").Value = 19431.666
$ws_ALC.Range("M137").Value = -3106.174199999999
$ws_ALC.Range("N137").Value = -24531.666

# Hunk 2: ARM row 74
$ws_ARM.Range("H74").Value = 2113.158
$ws_ARM.Range("I74").Value = 2268.1333
$ws_ARM.Range("K74").Value = 2268.1333
$ws_ARM.Range("M74").Value = -1394.1333

# Hunk 3: ARM row 77
$ws_ARM.Range("H77").Value = 2113.158
$ws_ARM.Range("I77").Value = 2268.1333
$ws_ARM.Range("K77").Value = 11340.6665
$ws_ARM.Range("M77").Value = -6972.666499999999

# Hunk 4: ARM row 101
$ws_ARM.Range("H101").Value = 20000
$ws_ARM.Range("J101").Value = 20000
$ws_ARM.Range("L101").Value = 20000
$ws_ARM.Range("N101").Value = -26490

# Hunk 5: BSM row 99
$ws_BSM.Range("H99").Value = 2339.1538
$ws_BSM.Range("I99").Value = 2178.889
$ws_BSM.Range("K99").Value = 2178.889
$ws_BSM.Range("M99").Value = -680.8890000000001

# Hunk 6: BSM row 104
$ws_BSM.Range("H104").Value = 0
$ws_BSM.Range("J104").Value = 0
$ws_BSM.Range("L104").Value = 0
$ws_BSM.Range("N104").ClearContents()

# Hunk 7: BSM row 107
$ws_BSM.Range("H107").Value = 1721.3684
$ws_BSM.Range("I107").Value = 1480
$ws_BSM.Range("J107").Value = 2053.25
$ws_BSM.Range("K107").Value = 1480
$ws_BSM.Range("L107").Value = 2053.25
$ws_BSM.Range("M107").Value = 440
$ws_BSM.Range("N107").Value = -5893.25

# Hunk 8: CRP row 58
$ws_CRP.Range("H58").Value = 3956122.8
$ws_CRP.Range("I58").Value = 7578193.5
$ws_CRP.Range("J58").Value = 4772.727
$ws_CRP.Range("K58").Value = 7578193.5
$ws_CRP.Range("L58").Value = 4772.727
$ws_CRP.Range("M58").Value = -7577990.5
$ws_CRP.Range("N58").Value = -5178.727

# Hunk 9: CRP row 99
$ws_CRP.Range("H99").Value = 2372
$ws_CRP.Range("I99").Value = 1486.7059
$ws_CRP.Range("J99").Value = 3626.1667
$ws_CRP.Range("K99").Value = 1486.7059
$ws_CRP.Range("L99").Value = 3626.1667
$ws_CRP.Range("M99").Value = 11.29410000000007
$ws_CRP.Range("N99").Value = -6622.1667

# Hunk 10: CRP row 107
$ws_CRP.Range("H107").Value = 731.5526
$ws_CRP.Range("I107").Value = 726.4231
$ws_CRP.Range("J107").Value = 742.6667
$ws_CRP.Range("K107").Value = 726.4231
$ws_CRP.Range("L107").Value = 742.6667
$ws_CRP.Range("M107").Value = 1193.5769
$ws_CRP.Range("N107").Value = -4582.6667

# Hunk 11: CRP row 126
$ws_CRP.Range("H126").Value = 2372
$ws_CRP.Range("I126").Value = 1486.7059
$ws_CRP.Range("J126").Value = 3626.1667
$ws_CRP.Range("K126").Value = 4460.1177
$ws_CRP.Range("L126").Value = 10878.5001
$ws_CRP.Range("M126").Value = -1990.1177
$ws_CRP.Range("N126").Value = -15818.5001

# Hunk 12: CRP row 136
$ws_CRP.Range("H136").Value = 3956122.8
$ws_CRP.Range("I136").Value = 7578193.5
$ws_CRP.Range("J136").Value = 4772.727
$ws_CRP.Range("K136").Value = 22734580.5
$ws_CRP.Range("L136").Value = 14318.181
$ws_CRP.Range("M136").Value = -22732030.5
$ws_CRP.Range("N136").Value = -19418.181

# Hunk 13: CUL row 5
$ws_CUL.Range("H5").Value = 3704250.8
$ws_CUL.Range("I5").Value = 357.91177
$ws_CUL.Range("K5").Value = 1073.73531
$ws_CUL.Range("M5").Value = -961.73531

# Hunk 14: CUL row 87
$ws_CUL.Range("H87").Value = 7252
$ws_CUL.Range("I87").Value = 2805.6
$ws_CUL.Range("J87").Value = 9722.223
$ws_CUL.Range("K87").Value = 8416.799999999999
$ws_CUL.Range("L87").Value = 29166.669
$ws_CUL.Range("M87").Value = -7168.799999999999
$ws_CUL.Range("N87").Value = -31662.669

# Hunk 15: CUL row 90
$ws_CUL.Range("H90").Value = 7252
$ws_CUL.Range("I90").Value = 2805.6
$ws_CUL.Range("J90").Value = 9722.223
$ws_CUL.Range("K90").Value = 25250.4
$ws_CUL.Range("L90").Value = 87500.007
$ws_CUL.Range("M90").Value = -19010.4
$ws_CUL.Range("N90").Value = -99980.007

# Hunk 16: CUL row 118
$ws_CUL.Range("H118").Value = 1888
$ws_CUL.Range("I118").Value = 1124
$ws_CUL.Range("J118").Value = 8000
$ws_CUL.Range("K118").Value = 3372
$ws_CUL.Range("L118").Value = 24000
$ws_CUL.Range("M118").Value = -2129
$ws_CUL.Range("N118").Value = -26486

# Hunk 17: CUL row 122
$ws_CUL.Range("H122").Value = 1001.6
$ws_CUL.Range("I122").Value = 457.5
$ws_CUL.Range("J122").Value = 1199.4546
$ws_CUL.Range("K122").Value = 4117.5
$ws_CUL.Range("L122").Value = 10795.0914
$ws_CUL.Range("M122").Value = -1667.5
$ws_CUL.Range("N122").Value = -15695.0914

# Hunk 18: CUL row 125
$ws_CUL.Range("H125").Value = 2706.1667
$ws_CUL.Range("I125").Value = 1666.6666
$ws_CUL.Range("J125").Value = 2854.6667
$ws_CUL.Range("K125").Value = 4999.9998
$ws_CUL.Range("L125").Value = 8564.000100000001
$ws_CUL.Range("M125").Value = -79.9997999999996
$ws_CUL.Range("N125").Value = -18404.0001

# Hunk 19: CUL row 134
$ws_CUL.Range("H134").Value = 3696.1538
$ws_CUL.Range("I134").Value = 3827.2
$ws_CUL.Range("J134").Value = 3517.4546
$ws_CUL.Range("K134").Value = 11481.6
$ws_CUL.Range("L134").Value = 10552.3638
$ws_CUL.Range("M134").Value = -6411.599999999999
$ws_CUL.Range("N134").Value = -20692.3638

# Hunk 20: CUL row 135
$ws_CUL.Range("H135").Value = 3704250.8
$ws_CUL.Range("I135").Value = 357.91177
$ws_CUL.Range("K135").Value = 3221.20593
$ws_CUL.Range("M135").Value = -686.2059300000001

# Hunk 21: CUL row 137
$ws_CUL.Range("H137").Value = 19254258
$ws_CUL.Range("I137").Value = 55556740
$ws_CUL.Range("J137").Value = 35297.41
$ws_CUL.Range("K137").Value = 166670220
$ws_CUL.Range("L137").Value = 105892.23
$ws_CUL.Range("M137").Value = -166665120
$ws_CUL.Range("N137").Value = -116092.23

# Hunk 22: GSM row 46
$ws_GSM.Range("H46").Value = 22800
$ws_GSM.Range("I46").Value = 5000
$ws_GSM.Range("J46").Value = 27250
$ws_GSM.Range("K46").Value = 5000
$ws_GSM.Range("L46").Value = 27250
$ws_GSM.Range("M46").Value = -4844
$ws_GSM.Range("N46").Value = -27562

# Hunk 23: GSM row 126
$ws_GSM.Range("H126").Value = 2868.2727
$ws_GSM.Range("I126").Value = 1860.2
$ws_GSM.Range("J126").Value = 3708.3333
$ws_GSM.Range("K126").Value = 5580.6
$ws_GSM.Range("L126").Value = 11124.9999
$ws_GSM.Range("M126").Value = -3110.6
$ws_GSM.Range("N126").Value = -16064.9999

# Hunk 24: LTW row 7
$ws_LTW.Range("H7").Value = 3904.7896
$ws_LTW.Range("I7").Value = 3484
$ws_LTW.Range("K7").Value = 3484
$ws_LTW.Range("M7").Value = -3372

# Hunk 25: LTW row 22
$ws_LTW.Range("H22").Value = 466.66666
$ws_LTW.Range("I22").Value = 300
$ws_LTW.Range("J22").Value = 550
$ws_LTW.Range("K22").Value = 300
$ws_LTW.Range("L22").Value = 550
$ws_LTW.Range("M22").Value = -5
$ws_LTW.Range("N22").Value = -1140

# Hunk 26: LTW row 27
$ws_LTW.Range("H27").Value = 466.66666
$ws_LTW.Range("I27").Value = 300
$ws_LTW.Range("J27").Value = 550
$ws_LTW.Range("K27").Value = 300
$ws_LTW.Range("L27").Value = 550
$ws_LTW.Range("M27").Value = -193
$ws_LTW.Range("N27").Value = -764

# Hunk 27: LTW row 40
$ws_LTW.Range("H40").Value = 3798.0356
$ws_LTW.Range("I40").Value = 3606.5908
$ws_LTW.Range("J40").Value = 4500
$ws_LTW.Range("K40").Value = 3606.5908
$ws_LTW.Range("L40").Value = 4500
$ws_LTW.Range("M40").Value = -3470.5908
$ws_LTW.Range("N40").Value = -4772

# Hunk 28: LTW row 126
$ws_LTW.Range("H126").Value = 3904.7896
$ws_LTW.Range("I126").Value = 3484
$ws_LTW.Range("K126").Value = 10452
$ws_LTW.Range("M126").Value = -7982

# Hunk 29: LTW row 132
$ws_LTW.Range("H132").Value = 4718.2163
$ws_LTW.Range("I132").Value = 4299.1113
$ws_LTW.Range("J132").Value = 5849.8
$ws_LTW.Range("K132").Value = 12897.3339
$ws_LTW.Range("L132").Value = 17549.4
$ws_LTW.Range("M132").Value = -10367.3339
$ws_LTW.Range("N132").Value = -22609.4

# Hunk 30: WVR row 95
$ws_WVR.Range("H95").Value = 47829.332
$ws_WVR.Range("J95").Value = 47829.332
$ws_WVR.Range("L95").Value = 47829.332
$ws_WVR.Range("N95").Value = -53321.332

# Hunk 31: WVR row 96
$ws_WVR.Range("H96").Value = 1999
$ws_WVR.Range("I96").Value = 1675.1666
$ws_WVR.Range("J96").Value = 2484.75
$ws_WVR.Range("K96").Value = 1675.1666
$ws_WVR.Range("L96").Value = 2484.75
$ws_WVR.Range("M96").Value = -302.1666
$ws_WVR.Range("N96").Value = -5230.75

# Hunk 32: WVR row 126
$ws_WVR.Range("H126").Value = 1311.9524
$ws_WVR.Range("I126").Value = 1074.2941
$ws_WVR.Range("J126").Value = 2322
$ws_WVR.Range("K126").Value = 3222.8823
$ws_WVR.Range("L126").Value = 6966
$ws_WVR.Range("M126").Value = -752.8823000000002
$ws_WVR.Range("N126").Value = -11906

# Hunk 33: WVR row 132
$ws_WVR.Range("H132").Value = 1751.5555
$ws_WVR.Range("I132").Value = 710.5625
$ws_WVR.Range("J132").Value = 3265.7273
$ws_WVR.Range("K132").Value = 2131.6875
$ws_WVR.Range("L132").Value = 9797.1819
$ws_WVR.Range("M132").Value = 398.3125
$ws_WVR.Range("N132").Value = -14857.1819

# Hunk 34: WVR row 140
$ws_WVR.Range("H140").Value = 59467.418
$ws_WVR.Range("J140").Value = 59467.418
$ws_WVR.Range("L140").Value = 59467.418
$ws_WVR.Range("N140").Value = -69827.41800000001

# Hunk 35: WVR row 141
$ws_WVR.Range("H141").Value = 50579.5
$ws_WVR.Range("J141").Value = 50579.5
$ws_WVR.Range("L141").Value = 50579.5
$ws_WVR.Range("N141").Value = -60939.5
